$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B21 value
$ws.Range("B21").Value = 20315

# Fill in row 22 with new data
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 15253
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
